$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hawaii row (row 12): city names had a stray "1" typo-suffix; fix them. ---
$ws.Range("B12").Value = "Honolulu"
$ws.Range("C12").Value = "East Honolulu"
$ws.Range("D12").Value = "Pearl City"
$ws.Range("E12").Value = "Hilo"
$ws.Range("F12").Value = "Kailua"

# --- Pennsylvania row (row 39): drop "Erie" (col E), pulling the last ---
# --- city, Harrisburg (col F), one column to the left, and clearing   ---
# --- the now-empty trailing cell.                                     ---
$ws.Range("E39").Value = "Harrisburg"
$ws.Range("F39").Value = $null

# --- Reflect where the user ended up looking/selecting afterwards. ---
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("E39").Select()
